$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the individual row-4 cells into row 5 so the new
# row matches the existing styling (currency format in column B, date
# format in column C) without disturbing column E.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new expense entry values
$ws.Range("A5").Value = "Amit Patankar"
$ws.Range("B5").Value = "`$10.90/year"
$ws.Range("C5").Value = Get-Date -Year 2014 -Month 1 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("D5").Value = "Business Cards (250 for advertising and clients)"

# Update the selected cell to reflect where the user would click next
$ws.Range("A6").Select() | Out-Null
